$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value (45180) for every data
# row (rows 2 through 407). Update it to the new serial value (45181)
# while leaving everything else (formatting, other columns) untouched.
$lastRow = 407
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
